$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.640.31'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.231.73'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.56%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.55'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.31%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.232.47'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.544'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.161'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.21%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.76'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.92%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.503'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000271'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '39.02'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.761.93'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.666.41'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.32%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.234.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.63%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.29'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.38%  '
$ws.Range("E19").Value = '  +1.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '507.97'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.88%  '
$ws.Range("E22").Value = '  -0.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.06'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.96%  '
$ws.Range("E24").Value = '  -2.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.88'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.162'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +78.85%  '
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("E28").Value = '  -0.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.09'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.36'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.37%  '
$ws.Range("B31").Value = 'Stacks'
$ws.Range("C31").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.93'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.49%  '
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '28.28'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.33%  '
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.16'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.42'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.10%  '
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '55.49'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.76%  '
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0790'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +14.48%  '
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '495.66'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.95%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.21'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0421'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.129'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.91%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.73'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.295'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.45%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.949.07'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.67%  '
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.46'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.43'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.26%  '
$ws.Range("E49").Value = '  +0.21%  '
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("E51").Value = '  -4.96%  '

Write-Output "applied"